# Revert responsive design implementation
# Restores rows that were previously filtered out of the SAG2 sensor data
# sheets: ROW50-FE-LIFTER (1), ROW50-MID-LIFTER (2), ROW11-FE-LIFTER (3),
# ROW11-MID-LIFTER (4).
#
# FE-LIFTER sheets (1 and 3): row 29's timestamp becomes a real date value
# (style "YYYY-MM-DD HH:MM:SS"), and three more rows (30-32) of sensor
# readings are appended; the last of those (row 32) keeps its timestamp as
# plain text, matching how row 29 originally looked before being converted.
#
# MID-LIFTER sheets (2 and 4): eighteen more rows (68-85) of sensor
# readings are appended after the existing last row (67).

$wb = $excel.ActiveWorkbook

$dateFormat = "YYYY-MM-DD HH:MM:SS"

# ---- FE-LIFTER sheets (1 and 3) ----------------------------------------

$feTimes = @(45729.58035023148, 45729.5803721875, 45729.58039546297)
$feRows  = @(29, 30, 31)

$feB = "0x01,0x90"
$feC = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
$feD = "0x01,0x90,"
$feE = "0x14"
$feF = 400
$feG = 568631262647113800000000.0
$feH = 400
$feI = 20

foreach ($sheetIndex in @(1, 3)) {
    $ws = $wb.Worksheets.Item($sheetIndex)

    for ($i = 0; $i -lt $feRows.Length; $i++) {
        $r = $feRows[$i]
        $ws.Range("A$r").Value = $feTimes[$i]
        $ws.Range("A$r").NumberFormat = $dateFormat
        $ws.Range("B$r").Value = $feB
        $ws.Range("C$r").Value = $feC
        $ws.Range("D$r").Value = $feD
        $ws.Range("E$r").Value = $feE
        $ws.Range("F$r").Value = $feF
        $ws.Range("G$r").Value = $feG
        $ws.Range("H$r").Value = $feH
        $ws.Range("I$r").Value = $feI
    }

    # Row 32 keeps its timestamp as plain text (unstyled), like row 29 did
    # before this edit.
    $ws.Range("A32").Value = "2025-03-14 01:55:46"
    $ws.Range("B32").Value = $feB
    $ws.Range("C32").Value = $feC
    $ws.Range("D32").Value = $feD
    $ws.Range("E32").Value = $feE
    $ws.Range("F32").Value = $feF
    $ws.Range("G32").Value = $feG
    $ws.Range("H32").Value = $feH
    $ws.Range("I32").Value = $feI
}

# ---- MID-LIFTER sheets (2 and 4) ---------------------------------------

$midTimes = @(
    45729.31518523148, 45729.31520722222, 45729.3152303588,
    45729.39866047454, 45729.39868246527, 45729.39870561343,
    45729.48213673611, 45729.48215891204, 45729.482181875,
    45729.56561206019, 45729.56563403935, 45729.56565724537,
    45729.64909206019, 45729.64911011574, 45729.64913337963,
    45729.73256368055, 45729.73258552083, 45729.73260887731
)

$midB = "0x01,0x90"
$midC = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$midD = "0x01,0x90,"
$midE = "0x19"
$midF = 400
$midG = 568631262647113800000000.0
$midH = 400
$midI = 25

foreach ($sheetIndex in @(2, 4)) {
    $ws = $wb.Worksheets.Item($sheetIndex)

    for ($i = 0; $i -lt $midTimes.Length; $i++) {
        $r = 68 + $i
        $ws.Range("A$r").Value = $midTimes[$i]
        $ws.Range("A$r").NumberFormat = $dateFormat
        $ws.Range("B$r").Value = $midB
        $ws.Range("C$r").Value = $midC
        $ws.Range("D$r").Value = $midD
        $ws.Range("E$r").Value = $midE
        $ws.Range("F$r").Value = $midF
        $ws.Range("G$r").Value = $midG
        $ws.Range("H$r").Value = $midH
        $ws.Range("I$r").Value = $midI
    }
}
